$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---
$ws.Range("A10").Value2 = 44566
$ws.Range("A10").NumberFormat = "d-mmm-yy"

$ws.Range("B10").Value2 = "Production"

$ws.Range("C10").Value2 = 134
$ws.Range("C10").HorizontalAlignment = $xlCenter
$ws.Range("C10").VerticalAlignment = $xlCenter

$ws.Range("D10").Value2 = 132
$ws.Range("D10").HorizontalAlignment = $xlCenter
$ws.Range("D10").VerticalAlignment = $xlCenter

$ws.Range("E10").Value2 = 2
$ws.Range("E10").HorizontalAlignment = $xlCenter
$ws.Range("E10").VerticalAlignment = $xlCenter

$ws.Range("F10").Value2 = "After execution all test cases pass"
$ws.Range("F10").WrapText = $true

$ws.Range("G10").Value2 = "Test cases initially fail because of page load affected by network"
$ws.Range("G10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 75

# --- Row 11 ---
$ws.Range("A11").Value2 = 44566
$ws.Range("A11").NumberFormat = "d-mmm-yy"

$ws.Range("B11").Value2 = "Development"

$ws.Range("C11").Value2 = 119
$ws.Range("C11").HorizontalAlignment = $xlCenter
$ws.Range("C11").VerticalAlignment = $xlCenter

$ws.Range("D11").Value2 = 119
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("D11").VerticalAlignment = $xlCenter

$ws.Range("E11").Value2 = 0
$ws.Range("E11").HorizontalAlignment = $xlCenter
$ws.Range("E11").VerticalAlignment = $xlCenter

# --- Selection / view state to match the final workbook state ---
$ws.Range("F11:G11").Select()
